$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start of 2025 author list / membership refresh: remove Ester Hammond's row
# entirely (row 3) so the sheet matches the JISCMail / collaboration
# database. Deleting the row shifts Philip Burrows and Manjit Dosanjh up one
# row each (rows 4 & 5 -> 3 & 4).
$ws.Rows(3).Delete() | Out-Null

# Amato Giaccia now also has a Stanford address on file.
$ws.Range("E2").Value = "amato.giaccia@oncology.ox.ac.uk; giaccia@stanford.edu"

# Manjit Dosanjh (now row 4) - normalise the two email addresses to be
# separated with a semicolon instead of the word "and".
$ws.Range("E4").Value = "Manjit.Dosanjh@cern.ch;  manjit.dosanjh@physics.ox.ac.uk"

# Leave the same selection behind that the interactive row-delete left.
$ws.Range("A3:XFD3").Select() | Out-Null
